# "reporte usuarios y report resumen"
#
# USUARIO ADMINISTRADOR sheet (sheet1): insert a new column before column C
# to make room for a "PERIODO" header column, update the title banner text
# (drop the hard-coded "MES DE ENERO 2024" in favour of the generic title
# already used on the USUARIO REGISTRADOR sheet), and move the selection.

$wb = $excel.ActiveWorkbook

$wsAdmin = $wb.Worksheets.Item("USUARIO ADMINISTRADOR")

# Insert a new column at C; this shifts the existing C:G columns (and their
# custom widths / merged title band / header row) one slot to the right,
# to D:H, exactly mirroring the diff's column layout change.
$wsAdmin.Columns.Item(3).Insert()

# New header cell for the inserted column.
$wsAdmin.Range("C4").Value = "PERIODO"

# The merged title band (now D2:G2) loses the literal month/year and keeps
# the generic report title - same wording used on the other sheet.
$wsAdmin.Range("D2").Value = "REPORTE DE CIERRE DE REGISTRO DE METAS FÍSICAS DE ACTIVIDADES OPERATIVAS "

# Move the active selection as recorded in the saved view state.
$wsAdmin.Range("C6").Select() | Out-Null
